$p = $ppt.ActivePresentation
$s = $p.Slides.Item(20)
$sh = $s.Shapes.Item("Content Placeholder 2")
$tf = $sh.TextFrame
$tr = $tf.TextRange

# --- Paragraph (lvl=1) "超越表格在查阅值右侧或另外一张工作表"
#     -> "查阅表格在查阅值右侧或另外一张工作表"
#     split into 3 runs: "查" | "阅表" | "格在查阅值右侧或另外一张工作表"
$para7 = $tr.Paragraphs(7,1)
$r1 = $para7.Characters(1,1)
$r1.Text = "查"

$para7b = $tr.Paragraphs(7,1)
$r2 = $para7b.Characters(2,2)
$r2.Text = "阅表"

# --- Paragraph (lvl=1) "按列排列数据" -> "按行排列数据"
#     split into 4 runs: "按行排" | "列" | "数" | "据"
$para8 = $tr.Paragraphs(8,1)
$q1 = $para8.Characters(1,3)
$q1.Text = "按行排"

$para8b = $tr.Paragraphs(8,1)
$q2 = $para8b.Characters(4,1)
$q2.Text = "列"

$para8c = $tr.Paragraphs(8,1)
$q3 = $para8c.Characters(5,1)
$q3.Text = "数"

$para8d = $tr.Paragraphs(8,1)
$q4 = $para8d.Characters(6,1)
$q4.Text = "据"

# --- Paragraph (lvl=1) "复制时绝对单元格引用" -> "复制时没有用绝对单元格引用"
#     split second run into 3: "制" | "时没有用绝" | "对单元格引用"
$para9 = $tr.Paragraphs(9,1)
$t1 = $para9.Characters(3,2)
$t1.Text = "时没有用绝"
